$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Repeated experiments with CPU_MHZ correctly set ---
# Columns C (start trigger @) and I (stop trigger @) were re-measured with
# the CPU_MHZ parameter corrected; downstream formulas (D/E and J/K/N/N26)
# recalculate automatically from these raw inputs.

$ws.Range("C2").Value = 3755997163
$ws.Range("I2").Value = 3755997163

$ws.Range("C3").Value = 4109757230
$ws.Range("I3").Value = 4110237245.7363024

$ws.Range("C4").Value = 7548146043
$ws.Range("I4").Value = 8070395988.4585695

$ws.Range("C5").Value = 5167907133
$ws.Range("I5").Value = 5167907133

$ws.Range("C6").Value = 2087828612
$ws.Range("I6").Value = 2090075548.276222

$ws.Range("C7").Value = 5446301669
$ws.Range("I7").Value = 5446301669

$ws.Range("C8").Value = 1705839944
$ws.Range("I8").Value = 1714972682.549557

$ws.Range("C9").Value = 2623150048
$ws.Range("I9").Value = 2712016433.0499949

$ws.Range("C10").Value = 3189419155
$ws.Range("I10").Value = 3189819169.0287313

$ws.Range("C11").Value = 4032763475
$ws.Range("I11").Value = 4032763475

$ws.Range("C12").Value = 3452605083
$ws.Range("I12").Value = 3604611157.9068036

$ws.Range("C13").Value = 2527010749
$ws.Range("I13").Value = 2527010749

$ws.Range("C14").Value = 3635111428
$ws.Range("I14").Value = 3799353904.9259934

$ws.Range("C15").Value = 2201504365
$ws.Range("I15").Value = 2201504365

$ws.Range("C16").Value = 2773880934
$ws.Range("I16").Value = 2784945277.1341786

$ws.Range("C17").Value = 2382629385
$ws.Range("I17").Value = 2383110690.8254185

$ws.Range("C18").Value = 2450028786
$ws.Range("I18").Value = 2450028786

$ws.Range("C19").Value = 4055687263
$ws.Range("I19").Value = 4056247441.8764958

$ws.Range("C20").Value = 1524569632
$ws.Range("I20").Value = 1524569632

$ws.Range("C21").Value = 1297593381
$ws.Range("I21").Value = 1773796385.7043729

$ws.Range("C22").Value = 1613351913
$ws.Range("I22").Value = 1613351913

$ws.Range("C23").Value = 1128858651
$ws.Range("I23").Value = 1706143232.1826692

# --- widen the new D/J columns to match their sibling measurement columns ---
$ws.Columns("D").ColumnWidth = $ws.Columns("B").ColumnWidth
$ws.Columns("J").ColumnWidth = $ws.Columns("H").ColumnWidth

# --- a block of extra (now-cleared) experiment rows below the summary,
#     left behind with a plain "0" number format on column F ---
$ws.Range("F28:F49").NumberFormat = "0"

# --- restore the cursor/selection position left by the editing session ---
$ws.Range("N36").Select() | Out-Null
